# The "books" sheet used to track a book's series name and its position
# within that series. Those two columns are being dropped from the sheet
# (collections are now inserted wholesale from Excel on start, so the
# series/numInSeries bookkeeping columns are no longer needed here).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("books")

# Columns C ("series") and D ("numInSeries") are removed entirely; the
# former column E ("isRead") shifts left to become column C.
$ws.Range("C1:D1").EntireColumn.Delete()
